$d = $word.ActiveDocument

# Locate the paragraph containing the CFG production for PROC_CALL.
# Its text is unique in the document, so match on the full trimmed text.
$targetText = "| PROC_CALL ID LEFT_PAREN CALL_PARAM RIGHT_PAREN"
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq $targetText) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph containing '$targetText'"
}

# Find the "CALL_PARAM" word inside that paragraph only (so we don't
# accidentally match "CALL_PARAMS" elsewhere in the document).
$callParamRange = $target.Range.Duplicate
$found = $callParamRange.Find.Execute("CALL_PARAM")
if (-not $found) {
    throw "Could not find 'CALL_PARAM' inside the target paragraph"
}

$s = $callParamRange.Start
$e = $callParamRange.End

# Insert the literal angle brackets around CALL_PARAM so the text reads
# "...LEFT_PAREN <CALL_PARAM> RIGHT_PAREN".
$callParamRange.InsertBefore("<")
$callParamRange.InsertAfter(">")

# After insertion, "<" occupies [s, s+1), "CALL_PARAM" occupies
# [s+1, e+1) and ">" occupies [e+1, e+2).
$ltRange = $d.Range($s, $s + 1)
$wordRange = $d.Range($s + 1, $e + 1)
$gtRange = $d.Range($e + 1, $e + 2)

foreach ($seg in @($ltRange, $wordRange, $gtRange)) {
    $seg.Font.Name = "Lucida Console"
    $seg.Font.Bold = 1
    $seg.Font.Color = 255
    $seg.Font.Size = 10
}
